$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.698.43"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "1.597.63"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'211.29"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").Value = "'0.245"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").Value = "'19.75"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "1.821.89"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.600.87"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").Value = "'65.22"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").Value = "26.698.95"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "'210.18"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "'6.74"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "'4.26"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Value = "'2.31"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "'146.56"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("D33").Value = "'0.665"
$ws.Range("E33").Value = "  -4.96%  "
$ws.Range("D35").Value = "1.298.25"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "'1.47"
$ws.Range("E37").Value = "  -4.79%  "
$ws.Range("D38").Value = "'0.0171"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("D39").Value = "'0.842"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").Value = "'63.66"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "1.733.79"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").Value = "'90.13"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").Value = "'0.873"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").Value = "'0.0988"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").Value = "'7.51"
$ws.Range("E51").Value = "  -0.40%  "
